$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data cells in row 3
$ws.Range("A3").Value = "倚天剑"
$ws.Range("B3").Value = 15647896523
$ws.Range("D3").Value = 15647896523
$ws.Range("G3").Value = "沈阳市8月底01分校"
$ws.Range("H3").Value = "学校"

# Update sheet view: move the active selection to G9
$ws.Range("G9").Select() | Out-Null

# Update workbook window size (bookViews windowHeight) to match the saved view
$win = $excel.ActiveWindow
$win.Height = 17655
$win.Width = 24045
